# Moved delivery organisation path to be before practitioner key on service contact.
#
# On the "Service Contacts" sheet, the "delivery_organisation_path" column
# (previously the last data column, R) is moved to sit right after
# "episode_key" (column C) and before "practitioner_key", i.e. it becomes
# the new column D. Every column that was between D and R shifts one to the
# right (E..R); the final "service_contact_tags" column (S) is unaffected.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Service Contacts")

# Make room for the relocated column right before the old column D
# (practitioner_key). This shifts practitioner_key..delivery_organisation_path
# (D..R) one column to the right, to E..S.
$ws.Columns("D").Insert() | Out-Null

# The delivery_organisation_path column (originally R) is now at S.
# Copy it into the freshly inserted column D (values + formatting)...
$ws.Range("S1:S3").Copy($ws.Range("D1:D3")) | Out-Null
$ws.Columns("D").ColumnWidth = $ws.Columns("S").ColumnWidth

# ...then remove the now-duplicate column at S, closing the gap so the
# trailing service_contact_tags column shifts back from T to S.
$ws.Columns("S").Delete() | Out-Null

# Reflect the column selection left behind by the move.
$ws.Columns("D").Select() | Out-Null

# A couple of other sheets were left with column F selected while this edit
# was being reviewed/cross-checked (both already have
# delivery_organisation_path in column F).
$wb.Worksheets.Item("K10+").Range("F1:F5").Select() | Out-Null
$wb.Worksheets.Item("K5").Range("F1:F5").Select() | Out-Null

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("Metadata").Activate() | Out-Null
